$d = $word.ActiveDocument

# --- 1. Remove the "VAT 0% - Art. 21 from the Bulgarian VAT Law" text.
#        Word tracks the last edit with the hidden "_GoBack" bookmark, so
#        re-anchor it onto the text's location (collapsing there once the
#        text is deleted) instead of leaving it on the trailing empty
#        paragraph. ---

$target = "VAT 0% - Art. 21 from the Bulgarian VAT Law"
$rng = $d.Content
$found = $rng.Find.Execute($target, $false, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks.Item("_GoBack").Delete()
    }
    $d.Bookmarks.Add("_GoBack", $rng)

    $rng.Text = ""
}

# --- 2. Mark a few styles as QuickStyles (w:qFormat), matching the
#        "Add to Quick Style list" toggle for Normal Table, List and
#        Table Grid. ---

for ($i = 1; $i -le $d.Styles.Count; $i++) {
    $s = $d.Styles.Item($i)
    $name = $s.NameLocal
    if ($name -eq "Normal Table" -or $name -eq "List" -or $name -eq "Table Grid") {
        $s.QuickStyle = $true
    }
}
